# "Multiply calculation succeeded! Congratulations!"
#
# Adds a new "Exponentiation" worksheet as the last tab (after "Minus"),
# fills it with the num/Index/Value table, applies the "Check Cell" style
# to the decimal-point cells (D4:D5), makes the new sheet the active one
# with selection on O12, and leaves the other sheets untouched.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip -------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Exponentiation"

# --- row 2: "num:" label, the number itself, and the strlen note -----------
$ws.Range("B2").Value = "num:"
$ws.Range("C2").Value = 35.612499999999997
$ws.Range("E2").Value = "strlen=7"

# --- row 4: index header row ------------------------------------------------
$ws.Range("A4").Value = "Index"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 6

# --- row 5: value row (D5 holds the decimal point placeholder) -------------
$ws.Range("A5").Value = "Value"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "."
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 5

# highlight the decimal-point column with the built-in "Check Cell" style
$ws.Range("D4:D5").Style = "检查单元格"

# visually separate the table with the thicker top/bottom border rows
$ws.Rows.Item(3).RowHeight = 14.4
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 14.4

# --- make the new sheet active, matching the saved selection/active tab ----
$ws.Activate()
[void]$ws.Range("O12").Select()
